# ui: move command execution responsibility from CommandBox to MainWindow
#
# The Ui class diagram's "execute" association arrow (and the small
# highlight rectangle sitting on top of the vertical UiPart connector
# line at the point the arrow crosses it) moves from the CommandBox row
# up to the MainWindow row, reflecting that MainWindow (not CommandBox)
# now calls Logic#execute(String).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$EMU_PER_PT = 12700

# Convert an EMU value to points the way PowerPoint's Left/Top/Width/Height
# properties expect. A tiny positive nudge is added so that the point value
# survives the host's internal round-trip back to EMU (single-precision,
# truncating) and still lands exactly on the requested EMU instead of being
# clipped one unit short.
function EmuToPt([double]$emu) {
    $pt = $emu / $EMU_PER_PT
    for ($nudge = 0; $nudge -lt 5000; $nudge++) {
        $candidate = $pt + $nudge * 0.0000001
        if ([int64]([single]$candidate * $EMU_PER_PT) -eq $emu) {
            return $candidate
        }
    }
    return $pt
}

# Locate the two shapes by their stable shape Ids (robust to z-order).
$highlight = $null
$connector = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -eq 143) { $highlight = $sh }
    if ($sh.Id -eq 116) { $connector = $sh }
}

# "Rectangle 142" - small highlight rectangle masking the line crossing.
# Reposition only (size unchanged): (5435896,2743200) -> (5422048,2339335) EMU.
$highlight.Left = EmuToPt 5422048
$highlight.Top  = EmuToPt 2339335

# "Freeform 115" - the diagonal/elbow connector line from the UiPart's
# owning class over to Logic. Reposition and resize:
# off (3687515,2828802) -> (3186477,2405681) EMU
# ext (3048000,203200)  -> (3537529,45719)   EMU
$connector.Left   = EmuToPt 3186477
$connector.Top    = EmuToPt 2405681
$connector.Width  = EmuToPt 3537529
$connector.Height = EmuToPt 45719
